# Automatische test-sync: 2025-06-20 13:30:50
# Append a new incoming-mail row to the "Logs" sheet and its matching
# category-count row to the "Dashboard" sheet, then extend the chart's
# source ranges so the new Dashboard row is included in the plot.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append row 13 (new mail about opening hours)
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A13").Value = "Wat zijn jullie openingstijden?"
$logs.Range("B13").Value = "mailmind.test@zohomail.eu"
$logs.Range("C13").Value = "Hallo, ik zou graag willen weten wat jullie openingstijden zijn. Dank je wel!"
$logs.Range("D13").Value = "Openingstijden / Locatie"
$logs.Range("E13").Value = "Beste klant,`nHartelijk dank voor uw vraag. Onze openingstijden zijn maandag t/m vrijdag van 09:00 tot 17:00 uur. Voor verdere informatie kunt u onze website raadplegen of telefonisch contact met ons opnemen.`nMet vriendelijke groet,`n[Naam van het bedrijf]"
$logs.Range("F13").Value = "2025-06-20 13:30:12"
$logs.Range("G13").Value = "Ja"

# Extend the "Categorie" / "Beantwoord" conditional formatting to cover
# the newly added row (was D2:D12 / G2:G12, now D2:D13 / G2:G13).
# Re-pointing a single rule's AppliesTo range re-points the whole
# <conditionalFormatting sqref="..."> block while leaving every
# cfRule (and its dxfId/priority) untouched.
$logs.Range("D2:D12").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D13"))
$logs.Range("G2:G12").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G13"))

# ---------------------------------------------------------------------
# 2. Dashboard sheet: append row 8 (count for the new category)
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A8").Value = "Openingstijden / Locatie"
$dash.Range("B8").Value = 1

# ---------------------------------------------------------------------
# 3. Chart: extend the series ranges from row 7 to row 8
# ---------------------------------------------------------------------
$chartObj = $dash.ChartObjects(1)
$series = $chartObj.Chart.SeriesCollection(1)
$series.XValues = "='Dashboard'!`$A`$2:`$A`$8"
$series.Values = "='Dashboard'!`$B`$2:`$B`$8"
